# repull data, push all data, mean calculation
# Update the dSF column (F) values for several rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 4
$ws.Range("F13").Value = 4
$ws.Range("F14").Value = 1
$ws.Range("F19").Value = -6
$ws.Range("F24").Value = -2
